$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 191.5
$ws.Range("I2").Value = 129.25
$ws.Range("J2").Value = 316
$ws.Range("K2").Value = 129.25
$ws.Range("L2").Value = 316
$ws.Range("M2").Value = -16.25
$ws.Range("N2").Value = -542

$ws.Range("H9").Value = 282.63635
$ws.Range("I9").Value = 257.375
$ws.Range("K9").Value = 257.375
$ws.Range("M9").Value = -88.375

$ws.Range("H100").Value = 2788.3635
$ws.Range("I100").Value = 2788.3635
$ws.Range("K100").Value = 2788.3635
$ws.Range("M100").Value = -2247.3635

$ws.Range("H132").Value = 2858.2222
$ws.Range("I132").Value = 3119.0454
$ws.Range("K132").Value = 9357.136200000001
$ws.Range("M132").Value = -6827.136200000001

$ws.Range("H137").Value = 3112.725
$ws.Range("I137").Value = 2106.5557
$ws.Range("J137").Value = 3935.9546
$ws.Range("K137").Value = 6319.6671
$ws.Range("L137").Value = 11807.8638
$ws.Range("M137").Value = -3769.6671
$ws.Range("N137").Value = -16907.8638

$ws.Range("H138").Value = 1858.51
$ws.Range("I138").Value = 941.1458
$ws.Range("J138").Value = 2705.3076
$ws.Range("K138").Value = 2823.4374
$ws.Range("L138").Value = 8115.9228
$ws.Range("M138").Value = 2316.5626
$ws.Range("N138").Value = -18395.9228

$ws.Range("H141").Value = 1410.2
$ws.Range("J141").Value = 3130.8
$ws.Range("L141").Value = 9392.400000000001
$ws.Range("N141").Value = -19752.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6616.347
$ws.Range("I32").Value = 5062.933
$ws.Range("K32").Value = 5062.933
$ws.Range("M32").Value = -4775.933

$ws.Range("H45").Value = 70229.13
$ws.Range("I45").Value = 92404
$ws.Range("J45").Value = 9248.25
$ws.Range("K45").Value = 92404
$ws.Range("L45").Value = 9248.25
$ws.Range("M45").Value = -92027
$ws.Range("N45").Value = -10002.25

$ws.Range("H61").Value = 3330.4285
$ws.Range("I61").Value = 2468.6667
$ws.Range("J61").Value = 4881.6
$ws.Range("K61").Value = 2468.6667
$ws.Range("L61").Value = 4881.6
$ws.Range("M61").Value = -2256.6667
$ws.Range("N61").Value = -5305.6

$ws.Range("H97").Value = 13112.956
$ws.Range("I97").Value = 10188.947
$ws.Range("K97").Value = 10188.947
$ws.Range("M97").Value = -9692.947

$ws.Range("H111").Value = 107879.5
$ws.Range("J111").Value = 107879.5
$ws.Range("L111").Value = 107879.5
$ws.Range("N111").Value = -116059.5

$ws.Range("H132").Value = 1550.4889
$ws.Range("I132").Value = 1165.317
$ws.Range("K132").Value = 3495.951
$ws.Range("M132").Value = -965.951

$ws.Range("H136").Value = 3330.4285
$ws.Range("I136").Value = 2468.6667
$ws.Range("J136").Value = 4881.6
$ws.Range("K136").Value = 7406.000100000001
$ws.Range("L136").Value = 14644.8
$ws.Range("M136").Value = -4856.000100000001
$ws.Range("N136").Value = -19744.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H32").Value = 28342
$ws.Range("I32").Value = 25026
$ws.Range("K32").Value = 25026
$ws.Range("M32").Value = -24642

$ws.Range("H64").Value = 2299.25
$ws.Range("J64").Value = 2500
$ws.Range("L64").Value = 2500
$ws.Range("N64").Value = -2950

$ws.Range("H67").Value = 2299.25
$ws.Range("J67").Value = 2500
$ws.Range("L67").Value = 2500
$ws.Range("N67").Value = -4060

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3154.975
$ws.Range("I31").Value = 2672.5789
$ws.Range("K31").Value = 2672.5789
$ws.Range("M31").Value = -2377.5789

$ws.Range("H34").Value = 3154.975
$ws.Range("I34").Value = 2672.5789
$ws.Range("K34").Value = 2672.5789
$ws.Range("M34").Value = -2470.5789

$ws.Range("H94").Value = 1210.3
$ws.Range("I94").Value = 807.6667
$ws.Range("K94").Value = 807.6667
$ws.Range("M94").Value = -356.6667

$ws.Range("H132").Value = 1957.8
$ws.Range("I132").Value = 1957.8
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5873.4
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3343.4
$ws.Range("N132").ClearContents()

$ws.Range("H134").Value = 32495.902
$ws.Range("I134").Value = 41781.13
$ws.Range("J134").Value = 5800.875
$ws.Range("K134").Value = 125343.39
$ws.Range("L134").Value = 17402.625
$ws.Range("M134").Value = -122808.39
$ws.Range("N134").Value = -22472.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 167782
$ws.Range("I46").Value = 417359
$ws.Range("J46").Value = 1397.3334
$ws.Range("K46").Value = 1252077
$ws.Range("L46").Value = 4192.0002
$ws.Range("M46").Value = -1251986
$ws.Range("N46").Value = -4374.0002

$ws.Range("H132").Value = 2888.9
$ws.Range("J132").Value = 2999.077
$ws.Range("L132").Value = 26991.693
$ws.Range("N132").Value = -32051.693

$ws.Range("H137").Value = 2628.7778
$ws.Range("I137").Value = 1771.5
$ws.Range("K137").Value = 5314.5
$ws.Range("M137").Value = -214.5

$ws.Range("H140").Value = 2108.6667
$ws.Range("I140").Value = 1622.375
$ws.Range("K140").Value = 4867.125
$ws.Range("M140").Value = 312.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 106283.336
$ws.Range("I122").Value = 144777.31
$ws.Range("J122").Value = 6199
$ws.Range("K122").Value = 434331.93
$ws.Range("L122").Value = 18597
$ws.Range("M122").Value = -431881.93
$ws.Range("N122").Value = -23497

$ws.Range("H132").Value = 2924.4146
$ws.Range("I132").Value = 2399.2188
$ws.Range("J132").Value = 4791.778
$ws.Range("K132").Value = 7197.6564
$ws.Range("L132").Value = 14375.334
$ws.Range("M132").Value = -4667.6564
$ws.Range("N132").Value = -19435.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 12714

$ws.Range("H43").Value = 7250
$ws.Range("I43").Value = 8500
$ws.Range("J43").Value = 6000
$ws.Range("K43").Value = 8500
$ws.Range("L43").Value = 6000
$ws.Range("M43").Value = -8307
$ws.Range("N43").Value = -6386

$ws.Range("H46").Value = 7032.4
$ws.Range("I46").Value = 6800
$ws.Range("J46").Value = 7264.8
$ws.Range("K46").Value = 6800
$ws.Range("L46").Value = 7264.8
$ws.Range("M46").Value = -6612
$ws.Range("N46").Value = -7640.8

$ws.Range("H132").Value = 6480
$ws.Range("I132").Value = 6013.075
$ws.Range("J132").Value = 8555.223
$ws.Range("K132").Value = 18039.225
$ws.Range("L132").Value = 25665.669
$ws.Range("M132").Value = -15509.225
$ws.Range("N132").Value = -30725.669

$ws.Range("H136").Value = 24877.666
$ws.Range("I136").Value = 34060.332
$ws.Range("J136").Value = 4675.8
$ws.Range("K136").Value = 102180.996
$ws.Range("L136").Value = 14027.4
$ws.Range("M136").Value = -99630.99600000001
$ws.Range("N136").Value = -19127.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 17000
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 17000
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 17000
$ws.Range("N18").Value = -17346
$ws.Range("M18").ClearContents()

$ws.Range("H25").Value = 21000
$ws.Range("J25").Value = 21000
$ws.Range("L25").Value = 21000
$ws.Range("N25").Value = -21586

$ws.Range("H96").Value = 1288.6666
$ws.Range("I96").Value = 947.25
$ws.Range("J96").Value = 1561.8
$ws.Range("K96").Value = 947.25
$ws.Range("L96").Value = 1561.8
$ws.Range("M96").Value = 425.75
$ws.Range("N96").Value = -4307.8

$ws.Range("H107").Value = 7927.9287
$ws.Range("I107").Value = 4363
$ws.Range("K107").Value = 13089
$ws.Range("M107").Value = -11169

$ws.Range("H126").Value = 1721.4482
$ws.Range("I126").Value = 1630.2916
$ws.Range("K126").Value = 4890.8748
$ws.Range("M126").Value = -2420.8748

$ws.Range("H132").Value = 6304.7334
$ws.Range("I132").Value = 6701.6924
$ws.Range("J132").Value = 3724.5
$ws.Range("K132").Value = 20105.0772
$ws.Range("L132").Value = 11173.5
$ws.Range("M132").Value = -17575.0772
$ws.Range("N132").Value = -16233.5
